# Update countries & provincias Spain
# Applies the daily data refresh: swaps six pairs of adjacent country rows
# (the country name moves, the newer country's numbers are refreshed while
# the other country's former numbers slide to the neighboring row), updates
# a handful of other country rows with refreshed numbers, and bumps the
# "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados ..." timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 20:45"

# --- Row 4: Estados Unidos (numbers refresh only) ---
$ws.Range("B4").Value = 6657598
$ws.Range("C4").Value = 21351
$ws.Range("D4").Value = 3928107
$ws.Range("E4").Value = 2531740
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 330
$ws.Range("H4").Value = 197751

# --- Row 5: India (numbers refresh only) ---
$ws.Range("B5").Value = 4750370
$ws.Range("C5").Value = 92991
$ws.Range("D5").Value = 3697905
$ws.Range("E5").Value = 973867
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1092
$ws.Range("H5").Value = 78598

# --- Rows 16/17: Reino Unido <-> Francia swap ---
$ws.Range("A16").Value = "Francia"
$ws.Range("B16").Value = 373911
$ws.Range("C16").Value = 10561
$ws.Range("D16").Value = 89059
$ws.Range("E16").Value = 253942
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 30910

$ws.Range("A17").Value = "Reino Unido"
$ws.Range("B17").Value = 365174
$ws.Range("C17").Value = 3497
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 41623

# --- Row 24: (numbers refresh only) ---
$ws.Range("B24").Value = 260286
$ws.Range("C24").Value = 561
$ws.Range("D24").Value = 234850
$ws.Range("E24").Value = 16012
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 9424

# --- Rows 60/61: Suiza <-> Uzbekistan swap ---
$ws.Range("A60").Value = "Uzbekistan"
$ws.Range("B60").Value = 46721
$ws.Range("C60").Value = 561
$ws.Range("D60").Value = 43359
$ws.Range("E60").Value = 2978
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 384

$ws.Range("A61").Value = "Suiza"
$ws.Range("B61").Value = 46704
$ws.Range("C61").Value = 465
$ws.Range("D61").Value = 38500
$ws.Range("E61").Value = 6184
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 2020

# --- Row 77: (numbers refresh only) ---
$ws.Range("D77").Value = 7936
$ws.Range("E77").Value = 15494

# --- Rows 79/80: Corea del Sur <-> Libia swap ---
$ws.Range("A79").Value = "Libia"
$ws.Range("B79").Value = 22348
$ws.Range("C79").Value = 440
$ws.Range("D79").Value = 2506
$ws.Range("E79").Value = 19488
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 354

$ws.Range("A80").Value = "Corea del Sur"
$ws.Range("B80").Value = 22055
$ws.Range("C80").Value = 136
$ws.Range("D80").Value = 18029
$ws.Range("E80").Value = 3671
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 5
$ws.Range("H80").Value = 355

# --- Rows 113/114: Mozambique <-> Suazilandia swap ---
$ws.Range("A113").Value = "Suazilandia"
$ws.Range("B113").Value = 5050
$ws.Range("C113").Value = 25
$ws.Range("D113").Value = 4188
$ws.Range("E113").Value = 764
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 98

$ws.Range("A114").Value = "Mozambique"
$ws.Range("B114").Value = 5040
$ws.Range("C114").Value = 122
$ws.Range("D114").Value = 2905
$ws.Range("E114").Value = 2100
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 4
$ws.Range("H114").Value = 35

# --- Rows 121/122: Cabo Verde <-> Cuba swap ---
$ws.Range("A121").Value = "Cuba"
$ws.Range("B121").Value = 4653
$ws.Range("C121").Value = 60
$ws.Range("D121").Value = 3878
$ws.Range("E121").Value = 667
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 2
$ws.Range("H121").Value = 108

$ws.Range("A122").Value = "Cabo Verde"
$ws.Range("B122").Value = 4651
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 4076
$ws.Range("E122").Value = 531
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 44
